$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 2 de Abril de 2020 a las 16:50'
$ws.Range("B4").Value = 216722
$ws.Range("C4").Value = 1719
$ws.Range("D4").Value = 8904
$ws.Range("E4").Value = 202678
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = 5140
$ws.Range("A7").Value = 'Alemania'
$ws.Range("B7").Value = 81728
$ws.Range("C7").Value = 3747
$ws.Range("D7").Value = 19175
$ws.Range("E7").Value = 61556
$ws.Range("F7").Value = 3936
$ws.Range("G7").Value = 66
$ws.Range("H7").Value = 997
$ws.Range("A8").Value = 'China'
$ws.Range("B8").Value = 81589
$ws.Range("C8").Value = 35
$ws.Range("D8").Value = 76408
$ws.Range("E8").Value = 1863
$ws.Range("F8").Value = 429
$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 3318
$ws.Range("B24").Value = 5119
$ws.Range("C24").Value = 242
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 5040
$ws.Range("B33").Value = 2692
$ws.Range("C33").Value = 138
$ws.Range("E33").Value = 2585
$ws.Range("G33").Value = 8
$ws.Range("H33").Value = 51
$ws.Range("A45").Value = 'Republica Dominicana'
$ws.Range("B45").Value = 1380
$ws.Range("C45").Value = 96
$ws.Range("D45").Value = 16
$ws.Range("E45").Value = 1304
$ws.Range("F45").Value = 147
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 60
$ws.Range("A46").Value = 'Mexico'
$ws.Range("B46").Value = 1378
$ws.Range("C46").Value = 163
$ws.Range("D46").Value = 35
$ws.Range("E46").Value = 1306
$ws.Range("F46").Value = 1
$ws.Range("G46").Value = 8
$ws.Range("H46").Value = 37
$ws.Range("A47").Value = 'Peru'
$ws.Range("B47").Value = 1323
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 394
$ws.Range("E47").Value = 882
$ws.Range("F47").Value = 49
$ws.Range("G47").Value = 9
$ws.Range("H47").Value = 47
$ws.Range("A48").Value = 'Islandia'
$ws.Range("B48").Value = 1319
$ws.Range("C48").Value = 99
$ws.Range("D48").Value = 284
$ws.Range("E48").Value = 1031
$ws.Range("F48").Value = 12
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 4
$ws.Range("A49").Value = 'Panama'
$ws.Range("B49").Value = 1317
$ws.Range("E49").Value = 1276
$ws.Range("F49").Value = 50
$ws.Range("H49").Value = 32
$ws.Range("B64").Value = 772
$ws.Range("C64").Value = 44
$ws.Range("D64").Value = 202
$ws.Range("E64").Value = 516
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 54
$ws.Range("A72").Value = 'Moldavia'
$ws.Range("B72").Value = 505
$ws.Range("C72").Value = 82
$ws.Range("D72").Value = 23
$ws.Range("E72").Value = 477
$ws.Range("F72").Value = 65
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 5
$ws.Range("A73").Value = 'Libano'
$ws.Range("B73").Value = 494
$ws.Range("C73").Value = 15
$ws.Range("D73").Value = 43
$ws.Range("E73").Value = 435
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 16
$ws.Range("A74").Value = 'Letonia'
$ws.Range("B74").Value = 458
$ws.Range("C74").Value = 12
$ws.Range("D74").Value = 1
$ws.Range("E74").Value = 457
$ws.Range("F74").Value = 3
$ws.Range("H74").Value = 0
$ws.Range("A75").Value = 'Bulgaria'
$ws.Range("B75").Value = 457
$ws.Range("C75").Value = 35
$ws.Range("D75").Value = 25
$ws.Range("E75").Value = 422
$ws.Range("F75").Value = 18
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 10
$ws.Range("A76").Value = 'Principado de Andorra'
$ws.Range("B76").Value = 428
$ws.Range("C76").Value = 38
$ws.Range("D76").Value = 10
$ws.Range("E76").Value = 403
$ws.Range("F76").Value = 12
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 15
$ws.Range("A77").Value = 'Eslovaquia'
$ws.Range("B77").Value = 426
$ws.Range("C77").Value = 26
$ws.Range("E77").Value = 420
$ws.Range("F77").Value = 3
$ws.Range("H77").Value = 1
$ws.Range("A78").Value = 'Tunez'
$ws.Range("D78").Value = 5
$ws.Range("E78").Value = 406
$ws.Range("F78").Value = 10
$ws.Range("H78").Value = 12
$ws.Range("A87").Value = 'Reunion'
$ws.Range("B87").Value = 308
$ws.Range("C87").Value = 27
$ws.Range("D87").Value = 40
$ws.Range("E87").Value = 268
$ws.Range("F87").Value = 3
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("A88").Value = 'Bielorrusia'
$ws.Range("B88").Value = 304
$ws.Range("C88").Value = 141
$ws.Range("D88").Value = 53
$ws.Range("E88").Value = 247
$ws.Range("F88").Value = 2
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 4
$ws.Range("A89").Value = 'Burkina Faso'
$ws.Range("B89").Value = 288
$ws.Range("C89").Value = 6
$ws.Range("D89").Value = 50
$ws.Range("E89").Value = 222
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 16
$ws.Range("A90").Value = 'Camerun'
$ws.Range("B90").Value = 284
$ws.Range("C90").Value = 51
$ws.Range("D90").Value = 10
$ws.Range("E90").Value = 267
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 7
$ws.Range("A99").Value = 'Ghana'
$ws.Range("B99").Value = 204
$ws.Range("C99").Value = 9
$ws.Range("D99").Value = 31
$ws.Range("E99").Value = 168
$ws.Range("H99").Value = 5
$ws.Range("A100").Value = 'Malta'
$ws.Range("B100").Value = 196
$ws.Range("C100").Value = 8
$ws.Range("D100").Value = 2
$ws.Range("E100").Value = 194
$ws.Range("F100").Value = 2
$ws.Range("H100").Value = 0
$ws.Range("B107").Value = 160
$ws.Range("C107").Value = 26
$ws.Range("E107").Value = 141
$ws.Range("A122").Value = 'Gibraltar'
$ws.Range("B122").Value = 88
$ws.Range("C122").Value = 7
$ws.Range("D122").Value = 46
$ws.Range("E122").Value = 42
$ws.Range("A123").Value = 'Ruanda'
$ws.Range("B123").Value = 82
$ws.Range("D123").Value = 0
$ws.Range("E123").Value = 82
$ws.Range("A162").Value = 'Siria'
$ws.Range("C162").Value = 6
$ws.Range("D162").Value = 0
$ws.Range("E162").Value = 14
$ws.Range("H162").Value = 2
$ws.Range("A163").Value = 'San Martin (Parte Holandesa)'
$ws.Range("B163").Value = 16
$ws.Range("D163").Value = 6
$ws.Range("E163").Value = 9
$ws.Range("H163").Value = 1
$ws.Range("A164").Value = 'Guinea Ecuatorial'
$ws.Range("B164").Value = 15
$ws.Range("D164").Value = 1
$ws.Range("E164").Value = 14
$ws.Range("A166").Value = 'Mongolia'
$ws.Range("B166").Value = 14
$ws.Range("D166").Value = 2
$ws.Range("A167").Value = 'Santa Lucia'
$ws.Range("A168").Value = 'Benin'
$ws.Range("B168").Value = 13
$ws.Range("D168").Value = 1
$ws.Range("A169").Value = 'Dominica'
$ws.Range("B169").Value = 12
$ws.Range("D169").Value = 0
$ws.Range("E169").Value = 12
$ws.Range("H169").Value = 0
$ws.Range("A170").Value = 'Curazao'
$ws.Range("B170").Value = 11
$ws.Range("D170").Value = 3
$ws.Range("E170").Value = 7
$ws.Range("H170").Value = 1
$ws.Range("A172").Value = 'Surinam'
$ws.Range("C172").Value = 0
$ws.Range("A174").Value = 'Seychelles'
$ws.Range("A175").Value = 'Granada'
$ws.Range("C175").Value = 1
$ws.Range("A176").Value = 'Mozambique'
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 10
$ws.Range("A177").Value = 'Groenlandia'
$ws.Range("D177").Value = 2
$ws.Range("H177").Value = 0
$ws.Range("A178").Value = 'Guinea-Bisau'
$ws.Range("A179").Value = 'Suazilandia'
